# The "reviews_count" column (column E) is being removed from the sheet.
# Deleting the whole column shifts every column to its right (F..K) one
# position to the left (F->E, G->F, H->G, I->H, J->I, K->J), which matches
# the target diff: reviews_average/latitude/longitude/is_permanently_closed/
# gmaps_link/latest_review_date all move one column left, and the sheet's
# dimension shrinks from A1:K41 to A1:J41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("E:E").Delete()
